$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in newly added testing-database utterances / rearranged cells ---
$ws.Range("F5").Value = 'flip 2 coins'
$ws.Range("F6").Value = 'oh yeah, flip a coin'
$ws.Range("O6").Value = 'kicking? I wanna do some kicking'
$ws.Range("F7").Value = 'do a coin flip'
$ws.Range("O7").Value = 'welcome to the salty spitoon, how tough are ya?'
$ws.Range("F8").Value = 'let''s flip a coin'
$ws.Range("G8").Value = 'roll a 12 sided die'
$ws.Range("J8").Value = 'define iridocyclitis'
$ws.Range("N8").Value = 'google september by earth wind and fire'
$ws.Range("O8").Value = 'I''m trying, but my cLEATS are stUCK in your coRNEAS'
$ws.Range("F9").Value = 'flip a coin right now'
$ws.Range("G9").Value = 'roll a d4 for me'
$ws.Range("J9").Value = 'what''s autodefenestration mean'
$ws.Range("N9").Value = 'search how to solve a rubik''s cube'
$ws.Range("O9").Value = 'you like krabby pattys don''t you squidward'
$ws.Range("B10").Value = 'it''s time for the time'
$ws.Range("H10").Value = 'set an alarm for 12PM'
$ws.Range("J10").Value = 'define duckie'
$ws.Range("N10").Value = 'look up where pineapple comes from'
$ws.Range("O10").Value = 'he''s just standing there…MENACINGLY'
$ws.Range("B11").Value = 'is it too late to go back to sleep?'
$ws.Range("I11").Value = 'kill the timer'
$ws.Range("N11").Value = 'google no time like the present'
$ws.Range("O11").Value = 'figure it out'
$ws.Range("M12").Value = '6 + 600'
$ws.Range("N12").Value = 'search for falling in reverse on google'
$ws.Range("O12").Value = 'ponder that for a moment'
$ws.Range("M13").Value = '8008 / 4'
$ws.Range("N13").Value = 'google when jac''o''lanterns were first made'
$ws.Range("O13").Value = 'kick the ball'
$ws.Range("M14").Value = '19 - 21'
$ws.Range("O14").Value = 'should we vote on it?'
$ws.Range("K15").Value = 'how''s my schedule look'
$ws.Range("O15").Value = 'Frankie!'
$ws.Range("K16").Value = 'how''s my schedule looking'
$ws.Range("K17").Value = 'how busy am I today'
$ws.Range("A20").Value = 'is it snowy today'
$ws.Range("E21").Value = 'am I broke or am I allowed to eat today'
$ws.Range("L21").Value = 'check my assignments for me'

# --- Column width tweaks (columns B and F got wider to fit new text) ---
$ws.Columns.Item(2).ColumnWidth = 27.498697916666668
$ws.Columns.Item(6).ColumnWidth = 16.830729166666668

# --- Update selection / view to the last-edited area ---
$ws.Range("O17").Select()
